# Update "paises" COVID country stats and re-sort by total cases.
# Commit message: "Update countries & provincias Spain"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "last updated" timestamp in A1 -------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 19:16"

# --- 2) Update the per-country numbers that changed ------------------------
# Each entry: Country name -> Casos totales, Nuevos casos, Casos activos,
#             Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @(
    @{ Name = "Estados Unidos";       B = 49906; C = 6172; D = 361; E = 48911; F = 1175; G = 81; H = 634 },
    @{ Name = "Noruega";              B = 2779;  C = 154;  D = 6;   E = 2761;  F = 44;   G = 2;  H = 12  },
    @{ Name = "Canada";               B = 2583;  C = 492;  D = 112; E = 2446;  F = 1;    G = 1;  H = 25  },
    @{ Name = "Japon";                B = 1193;  C = 65;   D = 285; E = 865;   F = 54;   G = 1;  H = 43  },
    @{ Name = "Pakistan";             B = 972;   C = 97;   D = 13;  E = 952;   F = 0;    G = 1;  H = 7   },
    @{ Name = "Polonia";              B = 844;   C = 95;   D = 1;   E = 834;   F = 3;    G = 1;  H = 9   },
    @{ Name = "India";                B = 536;   C = 37;   D = 40;  E = 486;   F = 0;    G = 0;  H = 10  },
    @{ Name = "Peru";                 B = 416;   C = 21;   D = 1;   E = 410;   F = 19;   G = 0;  H = 5   },
    @{ Name = "Argentina";            B = 301;   C = 0;    D = 51;  E = 244;   F = 0;    G = 2;  H = 6   },
    @{ Name = "Jordania";             B = 153;   C = 26;   D = 1;   E = 152;   F = 0;    G = 0;  H = 0   },
    @{ Name = "Republica de Chipre";  B = 124;   C = 8;    D = 3;   E = 118;   F = 3;    G = 2;  H = 3   },
    @{ Name = "Nigeria";              B = 44;    C = 4;    D = 2;   E = 41;    F = 0;    G = 0;  H = 1   },
    @{ Name = "Bolivia";              B = 29;    C = 2;    D = 0;   E = 29;    F = 0;    G = 0;  H = 0   },
    @{ Name = "Guatemala";            B = 21;    C = 1;    D = 0;   E = 20;    F = 0;    G = 0;  H = 1   },
    @{ Name = "Islas Caimanes";       B = 6;     C = 1;    D = 0;   E = 5;     F = 0;    G = 0;  H = 1   }
)

# Country names live in A4:A199 - look each one up and patch its row.
$searchRange = $ws.Range("A4:A199")
foreach ($u in $updates) {
    $cell = $searchRange.Find($u.Name)
    if ($cell -eq $null) {
        Write-Output "WARNING: country not found: $($u.Name)"
        continue
    }
    $r = $cell.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}

# --- 3) Re-sort the country table by "Casos totales" (column B), descending
$dataRange = $ws.Range("A4:H199")
$sortKey = $ws.Range("B4:B199")
$dataRange.Sort($sortKey, 2)

Write-Output "done"
